$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A header text was "hospital_id"; rename it to "id_hospital"
$ws.Range("A1").Value = "id_hospital"

# Move the selection/active cell to D13, matching the saved view state
$ws.Range("D13").Select()
